$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Profile" column with header and "Identity User" values
$ws.Range("G1").Value = "Profile"
$ws.Range("G2").Value = "Identity User"
$ws.Range("G3").Value = "Identity User"
$ws.Range("G4").Value = "Identity User"

# Build the custom font (JetBrains Mono, 9pt, greenish color) through a
# temporary named style so the font/name change is recorded as a single
# shared-string/style table entry instead of many incremental diffs.
$styleName = "IdentityUserFont"
$st = $wb.Styles.Add($styleName)
$st.Font.Name = "JetBrains Mono"

$rng = $ws.Range("G2:G4")
$rng.Style = $styleName
$wb.Styles.Item($styleName).Delete()

# Finish the font definition and vertically center the text.
$rng.Font.Family = 3
$rng.Font.Size = 9
$rng.Font.Color = 5867370
$rng.VerticalAlignment = -4108

$ws.Range("G4").Select()
